# Update "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 515a8dad... row
$wsOverview.Range("G3").Value = "2016-09-05 09:02:50"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 515a8dad... row
$wsZhCn.Range("H3").Value = "2016-09-05 09:02:45"
$wsZhCn.Range("K3").Value = "2016-09-05 09:03:22"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 515a8dad... row
$wsDeDe.Range("H3").Value = "2016-09-05 09:02:50"
$wsDeDe.Range("K3").Value = "2016-09-05 09:03:30"
